$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.830.00'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.209.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.83%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.610'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.70'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.43%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -2.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0905'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.95'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.541.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.49'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.208.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.780'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.796.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '228.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.46%  '
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '42.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.92%  '
$ws.Range("E28").Value = '  -3.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.08%  '
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0873'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.48%  '
$ws.Range("E34").Value = '  -1.71%  '
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0359'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.75%  '
$ws.Range("E37").Value = '  -2.60%  '
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("E40").Value = '  +18.96%  '
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.200'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("E43").Value = '  -4.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.35'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0979'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.466'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("E50").Value = '  -1.38%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.65'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.88%  '
